$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column S (19) slightly, as in the source commit
$ws.Columns.Item(19).ColumnWidth = 19.75

# Updated constants for removing 238U from 229Th and 230Th (columns N-S, rows 2-12)
$ws.Range("N2").Value = 1.319354151196508
$ws.Range("O2").Value = 0.5343093742747557
$ws.Range("P2").Value = 0.1699599420043786
$ws.Range("Q2").Value = 1.613124078274292
$ws.Range("R2").Value = 0.2255185977717201
$ws.Range("S2").Value = 1.726602086087542

$ws.Range("N3").Value = 0.03508119734104295
$ws.Range("O3").Value = 5.699227156799474
$ws.Range("P3").Value = 0.3294025880393198
$ws.Range("Q3").Value = 1.214208209894031
$ws.Range("R3").Value = 0.01151406902100858
$ws.Range("S3").Value = 2.442366103886794

$ws.Range("N4").Value = 1.319295276746676
$ws.Range("O4").Value = 0.700755786729371
$ws.Range("P4").Value = 0.1717532801460375
$ws.Range("Q4").Value = 1.289720296490281
$ws.Range("R4").Value = 0.2268636021468044
$ws.Range("S4").Value = 1.004357317888425

$ws.Range("N5").Value = 0.005800569700777472
$ws.Range("O5").Value = 2.870256697228719
$ws.Range("P5").Value = 0.225966131783047
$ws.Range("Q5").Value = 1.381881207195067
$ws.Range("R5").Value = 0.001249788646980554
$ws.Range("S5").Value = 3.023219177889351

$ws.Range("N6").Value = 1.313872396141101
$ws.Range("O6").Value = 0.5427985995815223
$ws.Range("P6").Value = 0.181231969044914
$ws.Range("Q6").Value = 1.277445604693938
$ws.Range("R6").Value = 0.2390518927769061
$ws.Range("S6").Value = 1.282867073244825

$ws.Range("N7").Value = 0.01772214606612332
$ws.Range("O7").Value = 1.007265011142688
$ws.Range("P7").Value = 0.249516127586646
$ws.Range("Q7").Value = 0.5610571919223687
$ws.Range("R7").Value = 0.004430472498451357
$ws.Range("S7").Value = 0.6991204502526103

$ws.Range("N8").Value = 1.31437756824995
$ws.Range("O8").Value = 0.4010354484850985
$ws.Range("P8").Value = 0.1790024332032209
$ws.Range("Q8").Value = 1.369577648543971
$ws.Range("R8").Value = 0.237144123290771
$ws.Range("S8").Value = 1.300904984569476

$ws.Range("N9").Value = 0.1435039070995674
$ws.Range("O9").Value = 0.3667519741005789
$ws.Range("P9").Value = 0.01853434975554646
$ws.Range("Q9").Value = 0.2678958047286344
$ws.Range("R9").Value = 0.002656456282872149
$ws.Range("S9").Value = 0.3014671536941368

$ws.Range("N10").Value = 1.313423515094221
$ws.Range("O10").Value = 0.6000345397267552
$ws.Range("P10").Value = 0.1778374152229192
$ws.Range("Q10").Value = 1.347700761114009
$ws.Range("R10").Value = 0.2338113376116541
$ws.Range("S10").Value = 1.249485075210555

$ws.Range("N11").Value = 0.1431288548060866
$ws.Range("O11").Value = 0.3371065744325605
$ws.Range("P11").Value = 0.0009627681596341573
$ws.Range("Q11").Value = 0.2554890408774075
$ws.Range("R11").Value = 0.0001382911205746749
$ws.Range("S11").Value = 0.2574720496270163

$ws.Range("N12").Value = 1.318565502507402
$ws.Range("O12").Value = 0.4544712691125686
$ws.Range("P12").Value = 0.1781347151424396
$ws.Range("Q12").Value = 1.295526126832855
$ws.Range("R12").Value = 0.2340996154520834
$ws.Range("S12").Value = 1.282918122977668

